# "fixed instruction image displays for the new training blocks"
#
# - Row 6 (full_task_wo_driving_training): enableSound yes -> no
# - Row 7 (full_task_training): instruction image ->
#       ./instructions_pilot/clutterlex_driving_training.png
# - Rows 8-10 (full_task_roboto / full_task_neuefrutigerworld / full_task_eurostile):
#       instruction image -> ./instructions_pilot/full_task_clutter.png
# - Selection moves to G7
# - Columns G and H widened to fit the longer instruction-image-file text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G6").Value = "no"

# write the "full_task_clutter.png" rows first so the new shared string for it
# gets allocated before the "clutterlex_driving_training.png" one (matches
# the shared-string ordering produced by the original edit)
$ws.Range("H8").Value = "./instructions_pilot/full_task_clutter.png"
$ws.Range("H9").Value = "./instructions_pilot/full_task_clutter.png"
$ws.Range("H10").Value = "./instructions_pilot/full_task_clutter.png"
$ws.Range("H7").Value = "./instructions_pilot/clutterlex_driving_training.png"

$ws.Range("G7").Select()

$ws.Columns.Item(7).ColumnWidth = 17.6
$ws.Columns.Item(8).ColumnWidth = 54.6
